# CTECH-2652 fix dates as strings (for cut labels)
#
# The holdings extract sheet is being re-shaped: a few columns that only
# ever held placeholder/duplicate data ("sub_holding_keys" = "{}",
# "SourcePortfolioId(default-Properties)", "SourcePortfolioScope(default-
# Properties)" and "cost_portfolio_ccy.amount" = 0) are dropped, and the
# remaining headers are renamed to the new camelCase extract schema:
#   instrument_uid                              -> luid
#   Name(default-Properties)                    -> instrumentName
#   holding_type                                 -> holdingType
#   units                                        -> units
#   settled_units                                -> settledUnits
#   cost.amount                                  -> costAmount
#   cost.currency                                -> costCurrency
#   cost_portfolio_ccy.currency                  -> portfolioCurrency
#
# The GBP cash-holding row's instrument name is also corrected from the
# (LUID-like) "CCY_GBP" to the real instrument name "GBP".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the unwanted columns. Delete from rightmost to leftmost so earlier
# deletions don't shift the column letters of ones still to be removed.
$ws.Range("L:L").Delete()   # cost_portfolio_ccy.amount (always 0)
$ws.Range("F:F").Delete()   # SourcePortfolioScope(default-Properties)
$ws.Range("E:E").Delete()   # SourcePortfolioId(default-Properties)
$ws.Range("C:C").Delete()   # sub_holding_keys (always "{}")

# Rename the remaining headers to the new schema.
$ws.Range("B1").Value = "luid"
$ws.Range("C1").Value = "instrumentName"
$ws.Range("D1").Value = "holdingType"
$ws.Range("E1").Value = "units"
$ws.Range("F1").Value = "settledUnits"
$ws.Range("G1").Value = "costAmount"
$ws.Range("H1").Value = "costCurrency"
$ws.Range("I1").Value = "portfolioCurrency"

# Correct the instrument name for the GBP cash holding row.
$ws.Range("C7").Value = "GBP"
